# Generate Report for Handback
# Updates the localization-status workbook to reflect a completed handback:
#  - Status changes from "Ready for handoff" to "Handed back: in sync with en-US"
#    on the Overview sheet and on each language sheet.
#  - Adds "Latest Target File" (E) and "Latest Handback File" (F) hyperlinked
#    entries for the two content rows on each language sheet.
#  - Stamps "Latest Handback DateTime" (G) with the handback timestamp.

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

$mdFileName  = "d0da193e-f564-4aa1-9caa-cd92e2e17c5e.md"
$mdUrl       = "https://github.com/OpenLocalizationTest/oltest/blob/a86268c3a6bf6552f4ab2d8552eae987023e9072/e2e/d0da193e-f564-4aa1-9caa-cd92e2e17c5e.md"

$xlfFileNameZh = "d0da193e-f564-4aa1-9caa-cd92e2e17c5e.8614f8fdde6b09b065ea7020fa16dfa3020e6125.zh-cn.xlf"
$xlfUrlZh      = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ae071574965e94fd36304c1dd9a573cd55af8319/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/d0da193e-f564-4aa1-9caa-cd92e2e17c5e.8614f8fdde6b09b065ea7020fa16dfa3020e6125.zh-cn.xlf"

$xlfFileNameDe = "d0da193e-f564-4aa1-9caa-cd92e2e17c5e.8614f8fdde6b09b065ea7020fa16dfa3020e6125.de-de.xlf"
$xlfUrlDe      = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/294eb8f1832a348eaecb588b5c1b4a1984f3ecb7/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/d0da193e-f564-4aa1-9caa-cd92e2e17c5e.8614f8fdde6b09b065ea7020fa16dfa3020e6125.de-de.xlf"

$handbackTimeZh = "2016-03-09 10:06:16"
$handbackTimeDe = "2016-03-09 10:06:36"

function Set-LinkCell($ws, $cellRef, $displayText, $url) {
    $ws.Hyperlinks.Add($ws.Range($cellRef), $url, [System.Type]::Missing, [System.Type]::Missing, $displayText) | Out-Null
    $ws.Range($cellRef).Font.Underline = 2
    $ws.Range($cellRef).Font.Color = 15570276
}

# ---- Overview sheet ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value2 = $statusText
$wsOverview.Range("C2").Value2 = $statusText
$wsOverview.Range("B3").Value2 = $statusText
$wsOverview.Range("C3").Value2 = $statusText

# ---- zh-cn sheet ----
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("B2").Value2 = $statusText
$wsZh.Range("B3").Value2 = $statusText

Set-LinkCell $wsZh "E2" $mdFileName $mdUrl
Set-LinkCell $wsZh "F2" $xlfFileNameZh $xlfUrlZh
Set-LinkCell $wsZh "E3" $mdFileName $mdUrl
Set-LinkCell $wsZh "F3" $xlfFileNameZh $xlfUrlZh

$wsZh.Range("G2").Value2 = $handbackTimeZh
$wsZh.Range("G3").Value2 = $handbackTimeZh

# ---- de-de sheet ----
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("B2").Value2 = $statusText
$wsDe.Range("B3").Value2 = $statusText

Set-LinkCell $wsDe "E2" $mdFileName $mdUrl
Set-LinkCell $wsDe "F2" $xlfFileNameDe $xlfUrlDe
Set-LinkCell $wsDe "E3" $mdFileName $mdUrl
Set-LinkCell $wsDe "F3" $xlfFileNameDe $xlfUrlDe

$wsDe.Range("G2").Value2 = $handbackTimeDe
$wsDe.Range("G3").Value2 = $handbackTimeDe

$wb.Save()
